$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.673.89'
$ws.Cells.Item(2, 5).Value = '  +0.10%  '
$ws.Cells.Item(3, 4).Value = '1.597.44'
$ws.Cells.Item(3, 5).Value = '  -0.12%  '
$ws.Cells.Item(4, 5).Value = '  +0.20%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '211.46'
$ws.Cells.Item(5, 5).Value = '  +0.28%  '
$ws.Cells.Item(6, 5).Value = '  -0.26%  '
$ws.Cells.Item(7, 5).Value = '  +0.10%  '
$ws.Cells.Item(8, 5).Value = '  +0.06%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.247'
$ws.Cells.Item(9, 5).Value = '  +0.64%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.49'
$ws.Cells.Item(10, 5).Value = '  -0.79%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0840'
$ws.Cells.Item(11, 5).Value = '  +0.36%  '
$ws.Cells.Item(12, 4).Value = '1.821.27'
$ws.Cells.Item(12, 5).Value = '  -0.12%  '
$ws.Cells.Item(13, 4).Value = '1.572.96'
$ws.Cells.Item(13, 5).Value = '  -1.57%  '
$ws.Cells.Item(14, 5).Value = '  +0.19%  '
$ws.Cells.Item(15, 5).Value = '  +0.45%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '65.10'
$ws.Cells.Item(16, 5).Value = '  +0.46%  '
$ws.Cells.Item(17, 4).Value = '26.640.15'
$ws.Cells.Item(17, 5).Value = '  -0.01%  '
$ws.Cells.Item(18, 4).Value = '0.0₃0737'
$ws.Cells.Item(18, 5).Value = '  +1.04%  '
$ws.Cells.Item(19, 5).Value = '  +0.18%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '209.15'
$ws.Cells.Item(20, 5).Value = '  -0.05%  '
$ws.Cells.Item(21, 5).Value = '  +4.56%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.29'
$ws.Cells.Item(22, 5).Value = '  +0.61%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '2.36'
$ws.Cells.Item(23, 5).Value = '  +2.66%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '143.88'
$ws.Cells.Item(25, 5).Value = '  -1.33%  '
$ws.Cells.Item(27, 5).Value = '  -1.71%  '
$ws.Cells.Item(28, 5).Value = '  -1.09%  '
$ws.Cells.Item(29, 5).Value = '  -0.02%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0513'
$ws.Cells.Item(30, 5).Value = '  +1.66%  '
$ws.Cells.Item(31, 5).Value = '  +0.15%  '
$ws.Cells.Item(32, 5).Value = '  -0.09%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.94'
$ws.Cells.Item(33, 5).Value = '  +0.57%  '
$ws.Cells.Item(34, 4).Value = '1.292.33'
$ws.Cells.Item(34, 5).Value = '  -0.31%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.620'
$ws.Cells.Item(35, 5).Value = '  -6.43%  '
$ws.Cells.Item(36, 5).Value = '  +0.53%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.49'
$ws.Cells.Item(37, 5).Value = '  -0.05%  '
$ws.Cells.Item(38, 5).Value = '  -0.34%  '
$ws.Cells.Item(39, 5).Value = '  -1.45%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.03'
$ws.Cells.Item(40, 5).Value = '  +15.38%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '5.46'
$ws.Cells.Item(41, 5).Value = '  +0.98%  '
$ws.Cells.Item(42, 5).Value = '  -0.44%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.782'
$ws.Cells.Item(43, 5).Value = '  -0.84%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '63.25'
$ws.Cells.Item(44, 5).Value = '  -0.94%  '
$ws.Cells.Item(45, 4).Value = '1.732.76'
$ws.Cells.Item(45, 5).Value = '  -0.19%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '91.11'
$ws.Cells.Item(46, 5).Value = '  +1.05%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.57'
$ws.Cells.Item(47, 5).Value = '  -2.97%  '
$ws.Cells.Item(48, 2).Value = 'Algorand'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.101'
$ws.Cells.Item(48, 5).Value = '  +0.83%  '
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0509'
$ws.Cells.Item(49, 5).Value = '  +0.99%  '
$ws.Cells.Item(50, 2).Value = 'USDD'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.00'
$ws.Cells.Item(50, 5).Value = '  +0.17%  '
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.38'
$ws.Cells.Item(51, 5).Value = '  -1.74%  '
